# Add two new vessel/project records (rows 13 and 14) to the "2025-1"
# sheet, matching the newly uploaded data: "E/P MODESTO 7" (A.S/0034) and
# "E/P TERESA" (A.S/0035).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: E/P MODESTO 7 -------------------------------------------
$ws.Range("A13").Value = "2025-1"
$ws.Range("D13").Value = "A.S/0034"
$ws.Range("E13").Value = "A.S/0034-125"
$ws.Range("B13").Value = "MODESTO 7"
$ws.Range("C13").Value = "Embarcación Pesquera"

# --- Row 14: E/P TERESA -----------------------------------------------
$ws.Range("A14").Value = "2025-1"
$ws.Range("B14").Value = "E/P TERESA"
$ws.Range("C14").Value = "Embarcación Pesquera"
$ws.Range("E14").Value = "A.S/0035-125"
$ws.Range("D14").Value = "A.S/0035"

# Widen column E (REDI) slightly so the longer "-125" codes fit.
$ws.Columns.Item(5).ColumnWidth = 12.6640625

# Leave the selection where the user's editing session ended up.
$ws.Range("G18").Select()
